$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# xlPasteFormats constant used throughout to copy an existing cell's exact
# visual style (so the engine reuses an existing cellXfs entry instead of
# minting a new one for every touched cell).
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Rows 11-15 : EARNED column (C) gets 1.25 for each existing SL row.
# These cells already carry style 13, so a plain value write is enough;
# the dependent "EARNED " helper column (G) recalculates automatically.
# ---------------------------------------------------------------------
$ws.Range("C11").Value = 1.25
$ws.Range("C12").Value = 1.25
$ws.Range("C13").Value = 1.25
$ws.Range("C14").Value = 1.25
$ws.Range("C15").Value = 1.25

# Row 16 : C16 currently carries the "last row" style (border on 3 sides
# only); once a value is entered it should match the interior style (13)
# used by the rest of the table body, so borrow that exact format first.
$ws.Range("C11").Copy() | Out-Null
$ws.Range("C16").PasteSpecial($xlPasteFormats)
$ws.Range("C16").Value = 1.25

# ---------------------------------------------------------------------
# Rows 17-18 : new dated entries. A17/A18 start out with the plain
# "0.000" numeric style; once a date is recorded they switch to the
# bordered mm/dd/yy style already used at A16 (style 41).
# ---------------------------------------------------------------------
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17").PasteSpecial($xlPasteFormats)
$ws.Range("A17").Value = 45254

$ws.Range("A16").Copy() | Out-Null
$ws.Range("A18").PasteSpecial($xlPasteFormats)
$ws.Range("A18").Value = 45284

$ws.Range("C11").Copy() | Out-Null
$ws.Range("C17").PasteSpecial($xlPasteFormats)
$ws.Range("C17").Value = 1.25

$ws.Range("C11").Copy() | Out-Null
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("C18").Value = 1.25

# ---------------------------------------------------------------------
# Row 19 : new "2024" year-separator row, mirroring the existing "2023"
# separator at A10 but using the date-column border (A16's border) and
# entered with a leading apostrophe so it is stored as literal text
# instead of being reinterpreted as the number 2024.
# ---------------------------------------------------------------------
$a19 = $ws.Range("A19")
$a16 = $ws.Range("A16")
$a19.Value = "'2024"
$a19.NumberFormat = $a16.NumberFormat
$a19.Font.Bold = $true
$a19.HorizontalAlignment = $a16.HorizontalAlignment
$a19.VerticalAlignment = $a16.VerticalAlignment
$a19.Borders.Item(7).LineStyle = $a16.Borders.Item(7).LineStyle
$a19.Borders.Item(8).LineStyle = $a16.Borders.Item(8).LineStyle
$a19.Borders.Item(8).Color = $a16.Borders.Item(8).Color
$a19.Borders.Item(10).LineStyle = $a16.Borders.Item(10).LineStyle
$a19.Borders.Item(10).Color = $a16.Borders.Item(10).Color

# ---------------------------------------------------------------------
# Row 20 : first 2024 SL entry - date, particulars, the 2-day "Absence
# Undertime  W/ Pay" offset, and the remarks note about the 01/04,05
# period.
# ---------------------------------------------------------------------
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A20").PasteSpecial($xlPasteFormats)
$ws.Range("A20").Value = 45315

$ws.Range("B20").Value = "SL(2-0-0)"
$ws.Range("H20").Value = 2
$ws.Range("K20").Value = "01/04,05/2024"

# ---------------------------------------------------------------------
# Table1 grows by one row: push the current last row (130, the
# special "bottom border" formatted row) down to 131 unchanged, then
# reset row 130 back to the plain interior-row style used throughout
# the table body (cloned from row 129) before the table is resized.
# ---------------------------------------------------------------------
$xlPasteAllExceptBorders = $xlPasteFormats
$earnedFormula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Placeholder write so row 131 physically exists before we paste into it.
$ws.Range("A131:K131").Value = "x"

$ws.Range("A130:K130").Copy() | Out-Null
$ws.Range("A131:K131").PasteSpecial($xlPasteFormats)
$ws.Range("A131:K131").ClearContents()
$ws.Range("G131").Formula = $earnedFormula

$ws.Range("A129:K129").Copy() | Out-Null
$ws.Range("A130:K130").PasteSpecial($xlPasteFormats)
$ws.Range("A130:K130").ClearContents()
$ws.Range("G130").Formula = $earnedFormula

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K131"))

$excel.CalculateFullRebuild()
